$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.465.01"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.971.94"
$ws.Range("E3").Value = "  -3.97%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.20"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  -4.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "56.89"
$ws.Range("E7").Value = "  +5.08%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.95"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.358"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0730"
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("E12").Value = "  -3.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.941"
$ws.Range("E13").Value = "  +5.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.14"
$ws.Range("E14").Value = "  -2.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.261.35"
$ws.Range("E15").Value = "  -3.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.23"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.971.34"
$ws.Range("E17").Value = "  -4.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.40"
$ws.Range("E18").Value = "  +4.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "35.413.33"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.09"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "232.21"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.12"
$ws.Range("E23").Value = "  -2.10%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("E25").Value = "  +19.75%  "
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.13"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.05"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.20"
$ws.Range("E29").Value = "  -3.97%  "
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.83"
$ws.Range("E31").Value = "  -3.80%  "
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0591"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0917"
$ws.Range("E34").Value = "  +10.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.26"
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.35"
$ws.Range("E36").Value = "  +8.38%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.13"
$ws.Range("E39").Value = "  +6.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.20"
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0211"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "91.10"
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.376.38"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.86"
$ws.Range("E46").Value = "  +1.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0880"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.49"
$ws.Range("E48").Value = "  +2.96%  "
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.25"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.46"
$ws.Range("E51").Value = "  +3.25%  "
